# Update the Uppmax/venue "linkoping" row (row 12) with the new LiU Campus
# address/location and refresh the selected cell + print page setup, as
# described in the commit "Updated uppmax project and venue address".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New venue address (lat/lon) for the Linköping location, row 12.
$ws.Range("D12").Value = 58.403739645581503
$ws.Range("E12").Value = 15.6223647575603

$newAddress = "LiU Campus Universitetssjukhuset`n58225 Linköping`nSweden`nMonday-Tuesday: [Room Papaver, Hus 511/001](https://link.mazemap.com/00mnumNU)`nWednesday-Thursday: [Room Dolomit, Hus 440](https://link.mazemap.com/up3GnjPm)`nFriday: [Room Antracit, Hus 440](https://link.mazemap.com/wXeFDYNR)"
$ws.Range("F12").Value = $newAddress
$ws.Range("F12").WrapText = $true

# Resize row 12 so the new multi-line address is fully visible.
$ws.Rows.Item(12).RowHeight = 119

# Move the active selection (as last left by the editing session).
[void]$ws.Range("A19").Select()

# Configure printing (page setup) for the sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
